$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 999
$ws.Range("I92").Value = 780.6667
$ws.Range("K92").Value = 780.6667
$ws.Range("M92").Value = 467.3333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3161.4
$ws.Range("I88").Value = 3333.3333
$ws.Range("J88").Value = 2903.5
$ws.Range("K88").Value = 3333.3333
$ws.Range("L88").Value = 2903.5
$ws.Range("M88").Value = -2927.3333
$ws.Range("N88").Value = -3715.5
$ws.Range("H91").Value = 3161.4
$ws.Range("I91").Value = 3333.3333
$ws.Range("J91").Value = 2903.5
$ws.Range("K91").Value = 3333.3333
$ws.Range("L91").Value = 2903.5
$ws.Range("M91").Value = -1929.3333
$ws.Range("N91").Value = -5711.5
$ws.Range("H132").Value = 1515.0227
$ws.Range("I132").Value = 1395.8055
$ws.Range("K132").Value = 4187.416499999999
$ws.Range("M132").Value = -1657.416499999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1899.2646
$ws.Range("I86").Value = 1591.0834
$ws.Range("J86").Value = 2638.9
$ws.Range("K86").Value = 1591.0834
$ws.Range("L86").Value = 2638.9
$ws.Range("M86").Value = -468.0834
$ws.Range("N86").Value = -4884.9
$ws.Range("H89").Value = 1899.2646
$ws.Range("I89").Value = 1591.0834
$ws.Range("J89").Value = 2638.9
$ws.Range("K89").Value = 7955.416999999999
$ws.Range("L89").Value = 13194.5
$ws.Range("M89").Value = -2339.416999999999
$ws.Range("N89").Value = -24426.5
$ws.Range("H99").Value = 1304.1
$ws.Range("I99").Value = 786.5
$ws.Range("J99").Value = 2080.5
$ws.Range("K99").Value = 786.5
$ws.Range("L99").Value = 2080.5
$ws.Range("M99").Value = 711.5
$ws.Range("N99").Value = -5076.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1495.28
$ws.Range("I31").Value = 1501.6471
$ws.Range("J31").Value = 1481.75
$ws.Range("K31").Value = 1501.6471
$ws.Range("L31").Value = 1481.75
$ws.Range("M31").Value = -1206.6471
$ws.Range("N31").Value = -2071.75
$ws.Range("H34").Value = 1495.28
$ws.Range("I34").Value = 1501.6471
$ws.Range("J34").Value = 1481.75
$ws.Range("K34").Value = 1501.6471
$ws.Range("L34").Value = 1481.75
$ws.Range("M34").Value = -1299.6471
$ws.Range("N34").Value = -1885.75
$ws.Range("H62").Value = 2866.6667
$ws.Range("J62").Value = 2800
$ws.Range("L62").Value = 2800
$ws.Range("N62").Value = -4048
$ws.Range("H65").Value = 2866.6667
$ws.Range("J65").Value = 2800
$ws.Range("L65").Value = 14000
$ws.Range("N65").Value = -20240
$ws.Range("H86").Value = 23812136
$ws.Range("I86").Value = 1920.25
$ws.Range("J86").Value = 38464576
$ws.Range("K86").Value = 1920.25
$ws.Range("L86").Value = 38464576
$ws.Range("M86").Value = -797.25
$ws.Range("N86").Value = -38466822
$ws.Range("H89").Value = 23812136
$ws.Range("I89").Value = 1920.25
$ws.Range("J89").Value = 38464576
$ws.Range("K89").Value = 9601.25
$ws.Range("L89").Value = 192322880
$ws.Range("M89").Value = -3985.25
$ws.Range("N89").Value = -192334112

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 100000760
$ws.Range("J34").Value = 166667700
$ws.Range("L34").Value = 500003100
$ws.Range("N34").Value = -500003268
$ws.Range("H92").Value = 1035.8182
$ws.Range("I92").Value = 450.2
$ws.Range("J92").Value = 1523.8334
$ws.Range("K92").Value = 1350.6
$ws.Range("L92").Value = 4571.5002
$ws.Range("M92").Value = -102.5999999999999
$ws.Range("N92").Value = -7067.5002
$ws.Range("H93").Value = 21880.5
$ws.Range("I93").Value = 40512
$ws.Range("J93").Value = 3249
$ws.Range("K93").Value = 121536
$ws.Range("L93").Value = 9747
$ws.Range("M93").Value = -119664
$ws.Range("N93").Value = -13491
$ws.Range("H96").Value = 70707200
$ws.Range("J96").Value = 70707200
$ws.Range("L96").Value = 212121600
$ws.Range("N96").Value = -212125718
$ws.Range("H113").Value = 503.65714
$ws.Range("I113").Value = 498.09525
$ws.Range("J113").Value = 512
$ws.Range("K113").Value = 1494.28575
$ws.Range("L113").Value = 1536
$ws.Range("M113").Value = 675.71425
$ws.Range("N113").Value = -5876
$ws.Range("H131").Value = 2002859.1
$ws.Range("J131").Value = 2129621.2
$ws.Range("L131").Value = 6388863.600000001
$ws.Range("N131").Value = -6398943.600000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4731.846
$ws.Range("I70").Value = 4931.5
$ws.Range("J70").Value = 4066.3333
$ws.Range("K70").Value = 4931.5
$ws.Range("L70").Value = 4066.3333
$ws.Range("M70").Value = -4661.5
$ws.Range("N70").Value = -4606.3333
$ws.Range("H73").Value = 4731.846
$ws.Range("I73").Value = 4931.5
$ws.Range("J73").Value = 4066.3333
$ws.Range("K73").Value = 4931.5
$ws.Range("L73").Value = 4066.3333
$ws.Range("M73").Value = -3995.5
$ws.Range("N73").Value = -5938.3333
$ws.Range("H80").Value = 12263.637
$ws.Range("I80").Value = 2414.2856
$ws.Range("J80").Value = 29500
$ws.Range("K80").Value = 2414.2856
$ws.Range("L80").Value = 29500
$ws.Range("M80").Value = -1416.2856
$ws.Range("N80").Value = -31496
$ws.Range("H83").Value = 12263.637
$ws.Range("I83").Value = 2414.2856
$ws.Range("J83").Value = 29500
$ws.Range("K83").Value = 12071.428
$ws.Range("L83").Value = 147500
$ws.Range("M83").Value = -7079.428
$ws.Range("N83").Value = -157484
$ws.Range("H97").Value = 1944.7333
$ws.Range("I97").Value = 1690
$ws.Range("J97").Value = 2326.8333
$ws.Range("K97").Value = 1690
$ws.Range("L97").Value = 2326.8333
$ws.Range("M97").Value = -1194
$ws.Range("N97").Value = -3318.8333
$ws.Range("H132").Value = 1851.6957
$ws.Range("I132").Value = 1422.0625
$ws.Range("J132").Value = 2833.7144
$ws.Range("K132").Value = 4266.1875
$ws.Range("L132").Value = 8501.143199999999
$ws.Range("M132").Value = -1736.1875
$ws.Range("N132").Value = -13561.1432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1283.3334
$ws.Range("I46").Value = 1352.8235
$ws.Range("J46").Value = 1114.5714
$ws.Range("K46").Value = 1352.8235
$ws.Range("L46").Value = 1114.5714
$ws.Range("M46").Value = -1164.8235
$ws.Range("N46").Value = -1490.5714
$ws.Range("H55").Value = 400.16666
$ws.Range("I55").Value = 467
$ws.Range("J55").Value = 333.33334
$ws.Range("K55").Value = 467
$ws.Range("L55").Value = 333.33334
$ws.Range("M55").Value = -294
$ws.Range("N55").Value = -679.33334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1470.659
$ws.Range("I122").Value = 1116.3334
$ws.Range("J122").Value = 2533.6365
$ws.Range("K122").Value = 3349.0002
$ws.Range("L122").Value = 7600.9095
$ws.Range("M122").Value = -899.0001999999999
$ws.Range("N122").Value = -12500.9095
$ws.Range("H126").Value = 2549.5454
$ws.Range("I126").Value = 3182.5
$ws.Range("J126").Value = 861.6667
$ws.Range("K126").Value = 9547.5
$ws.Range("L126").Value = 2585.0001
$ws.Range("M126").Value = -7077.5
$ws.Range("N126").Value = -7525.0001
